# Schedulator/Programs.xlsx — "Added program management to program director"
#
# 1. Rename the existing sheet and add a second sheet for the January
#    offering of the program (kept empty, as in the target workbook).
# 2. On the (now renamed) first sheet, insert a new "Septemeber" column
#    before the existing numeric value in row 1 (shifting it from C1 to D1).
# 3. Restore the original sheet/selection as the active one, with the
#    selection moved to D2 (matching the edited workbook).

$wb = $excel.ActiveWorkbook

# --- Sheets -----------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "SOEN COMP GAMES SEPT"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "SOEN COMP GAMES JAN"

# --- Row 1 on the first sheet: insert "Septemeber" in column C, ------
# --- pushing the old C1 value (120) into D1. --------------------------
$ws1.Cells.Item(1, 4).Value = $ws1.Cells.Item(1, 3).Value()
$ws1.Cells.Item(1, 3).Value = "Septemeber"

# --- Restore sheet1 as the active sheet/selection ----------------------
$ws1.Activate()
$null = $ws1.Select()
$null = $ws1.Range("D2").Select()
